$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 12 already exists but gets new/updated values across columns A:K.
# Rows 13-15 are brand new rows with the same "Alexnet / All 32" pattern
# but differing Cache Size (col J) and Hit Rate (col K).
# ---------------------------------------------------------------------

$newRows = @(
    @{ Row = 12; B = 0.502; C = 0.61; D = 0.012; J = 1000;  K = "10-50%"  },
    @{ Row = 13; B = 0.502; C = 0.61; D = 0.012; J = 2000;  K = "40-80%"  },
    @{ Row = 14; B = 0.502; C = 0.61; D = 0.012; J = 5000;  K = "60-95%"  },
    @{ Row = 15; B = 0.502; C = 0.61; D = 0.012; J = 10000; K = "90-100%" }
)

foreach ($item in $newRows) {
    $r = $item.Row

    if ($r -ne 12) {
        # New rows inherit the same row height as row 12 (13.8pt).
        $ws.Rows.Item($r).RowHeight = 13.8
        $ws.Cells.Item($r, 1).Value = "Alexnet"
    }

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 2).NumberFormat = "0.00%"

    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 3).NumberFormat = "0.00%"

    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 4).NumberFormat = "0.00%"

    $ws.Cells.Item($r, 5).Value = "All 32"

    $ws.Cells.Item($r, 6).Value = 0.015
    $ws.Cells.Item($r, 6).NumberFormat = "0.00%"

    $ws.Cells.Item($r, 7).Value = 0.0001
    $ws.Cells.Item($r, 7).NumberFormat = "0.00E+00"

    $ws.Cells.Item($r, 8).Value = 0.015
    $ws.Cells.Item($r, 8).NumberFormat = "0.00%"

    $ws.Cells.Item($r, 9).Value = 16

    $ws.Cells.Item($r, 10).Value = $item.J

    $ws.Cells.Item($r, 11).Value = $item.K
}

# Move the active selection, matching the workbook's last saved cursor
# position after the edits were made.
$ws.Range("F16").Select()
